$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1922.6666
$ws.Range("I4").Value = 1922.6666
$ws.Range("K4").Value = 1922.6666
$ws.Range("M4").Value = -1808.6666
$ws.Range("H17").Value = 1575
$ws.Range("J17").Value = 1500
$ws.Range("L17").Value = 4500
$ws.Range("N17").Value = -4836
$ws.Range("H18").Value = 19144.4
$ws.Range("J18").Value = 2922
$ws.Range("L18").Value = 2922
$ws.Range("N18").Value = -3490
$ws.Range("H28").Value = 1688.75
$ws.Range("I28").Value = 1688.75
$ws.Range("K28").Value = 1688.75
$ws.Range("M28").Value = -1203.75
$ws.Range("H33").Value = 214.72728
$ws.Range("I33").Value = 248
$ws.Range("J33").Value = 65
$ws.Range("K33").Value = 248
$ws.Range("L33").Value = 65
$ws.Range("M33").Value = -19
$ws.Range("N33").Value = -523
$ws.Range("H43").Value = 2738.3333
$ws.Range("J43").Value = 2738.3333
$ws.Range("L43").Value = 2738.3333
$ws.Range("N43").Value = -2876.3333
$ws.Range("H44").Value = 19499
$ws.Range("J44").Value = 19499
$ws.Range("L44").Value = 19499
$ws.Range("N44").Value = -20423
$ws.Range("H53").Value = 137.33333
$ws.Range("I53").Value = 76.833336
$ws.Range("K53").Value = 76.833336
$ws.Range("M53").Value = 560.166664
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("N79").ClearContents()
$ws.Range("H88").Value = 2635.5715
$ws.Range("I88").Value = 3287.25
$ws.Range("J88").Value = 1766.6666
$ws.Range("K88").Value = 3287.25
$ws.Range("L88").Value = 1766.6666
$ws.Range("M88").Value = -2881.25
$ws.Range("N88").Value = -2578.6666
$ws.Range("H91").Value = 2635.5715
$ws.Range("I91").Value = 3287.25
$ws.Range("J91").Value = 1766.6666
$ws.Range("K91").Value = 3287.25
$ws.Range("L91").Value = 1766.6666
$ws.Range("M91").Value = -1883.25
$ws.Range("N91").Value = -4574.6666
$ws.Range("H98").Value = 1500
$ws.Range("I98").Value = 1500
$ws.Range("K98").Value = 1500
$ws.Range("M98").Value = -2
$ws.Range("H111").Value = 395
$ws.Range("I111").Value = 290
$ws.Range("J111").Value = 500
$ws.Range("K111").Value = 870
$ws.Range("L111").Value = 1500
$ws.Range("M111").Value = 2197
$ws.Range("N111").Value = -7634
$ws.Range("H122").Value = 1500
$ws.Range("I122").Value = 1500
$ws.Range("K122").Value = 4500
$ws.Range("M122").Value = -2050
$ws.Range("H132").Value = 4303.8667
$ws.Range("J132").Value = 6999.8
$ws.Range("L132").Value = 20999.4
$ws.Range("N132").Value = -26059.4
$ws.Range("H135").Value = 2234.375
$ws.Range("I135").Value = 2247.8333
$ws.Range("K135").Value = 20230.4997
$ws.Range("M135").Value = -17695.4997

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2466.111
$ws.Range("I2").Value = 2782.8333
$ws.Range("K2").Value = 2782.8333
$ws.Range("M2").Value = -2669.8333
$ws.Range("H61").Value = 3386.7778
$ws.Range("I61").Value = 3386.7778
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3386.7778
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -3174.7778
$ws.Range("N61").ClearContents()
$ws.Range("H74").Value = 2638.5
$ws.Range("I74").Value = 2515.4285
$ws.Range("J74").Value = 3500
$ws.Range("K74").Value = 2515.4285
$ws.Range("L74").Value = 3500
$ws.Range("M74").Value = -1641.4285
$ws.Range("N74").Value = -5248
$ws.Range("H77").Value = 2638.5
$ws.Range("I77").Value = 2515.4285
$ws.Range("J77").Value = 3500
$ws.Range("K77").Value = 12577.1425
$ws.Range("L77").Value = 17500
$ws.Range("M77").Value = -8209.1425
$ws.Range("N77").Value = -26236
$ws.Range("H110").Value = 1763.4445
$ws.Range("I110").Value = 1132.8
$ws.Range("K110").Value = 1132.8
$ws.Range("M110").Value = 912.2
$ws.Range("H116").Value = 2466.111
$ws.Range("I116").Value = 2782.8333
$ws.Range("K116").Value = 2782.8333
$ws.Range("M116").Value = -488.8332999999998
$ws.Range("H132").Value = 3176.8
$ws.Range("I132").Value = 2971.25
$ws.Range("K132").Value = 8913.75
$ws.Range("M132").Value = -6383.75
$ws.Range("H136").Value = 3386.7778
$ws.Range("I136").Value = 3386.7778
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 10160.3334
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -7610.3334
$ws.Range("N136").ClearContents()

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2466.111
$ws.Range("I3").Value = 2782.8333
$ws.Range("K3").Value = 2782.8333
$ws.Range("M3").Value = -2668.8333
$ws.Range("H88").Value = 33350
$ws.Range("J88").Value = 42800
$ws.Range("L88").Value = 42800
$ws.Range("N88").Value = -43612
$ws.Range("H91").Value = 33350
$ws.Range("J91").Value = 42800
$ws.Range("L91").Value = 42800
$ws.Range("N91").Value = -45608
$ws.Range("H134").Value = 4824.75
$ws.Range("I134").Value = 3356.2856
$ws.Range("J134").Value = 8251.166999999999
$ws.Range("K134").Value = 10068.8568
$ws.Range("L134").Value = 24753.501
$ws.Range("M134").Value = -7533.856800000001
$ws.Range("N134").Value = -29823.501

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10000.3
$ws.Range("I31").Value = 11002
$ws.Range("K31").Value = 11002
$ws.Range("M31").Value = -10707
$ws.Range("H34").Value = 10000.3
$ws.Range("I34").Value = 11002
$ws.Range("K34").Value = 11002
$ws.Range("M34").Value = -10800
$ws.Range("H105").Value = 1755
$ws.Range("I105").Value = 1755
$ws.Range("K105").Value = 1755
$ws.Range("M105").Value = -8
$ws.Range("H122").Value = 807.1111
$ws.Range("I122").Value = 773.5714
$ws.Range("K122").Value = 2320.7142
$ws.Range("M122").Value = 129.2857999999997
$ws.Range("H132").Value = 2631
$ws.Range("I132").Value = 2631
$ws.Range("K132").Value = 7893
$ws.Range("M132").Value = -5363

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 1954.6666
$ws.Range("I137").Value = 924
$ws.Range("J137").Value = 2779.2
$ws.Range("K137").Value = 2772
$ws.Range("L137").Value = 8337.599999999999
$ws.Range("M137").Value = 2328
$ws.Range("N137").Value = -18537.6

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 110
$ws.Range("I2").Value = 128
$ws.Range("K2").Value = 128
$ws.Range("M2").Value = -15
$ws.Range("H5").Value = 1669.3334
$ws.Range("I5").Value = 1669.3334
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1669.3334
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -1557.3334
$ws.Range("N5").ClearContents()
$ws.Range("H102").Value = 2478.0715
$ws.Range("I102").Value = 2108.1667
$ws.Range("K102").Value = 2108.1667
$ws.Range("M102").Value = -486.1667000000002
$ws.Range("H126").Value = 15000
$ws.Range("J126").Value = 15000
$ws.Range("L126").Value = 45000
$ws.Range("N126").Value = -49940

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4981
$ws.Range("I7").Value = 4701.1113
$ws.Range("K7").Value = 4701.1113
$ws.Range("M7").Value = -4589.1113
$ws.Range("H16").Value = 8999.5
$ws.Range("I16").Value = 8999.5
$ws.Range("K16").Value = 8999.5
$ws.Range("M16").Value = -8829.5
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 750
$ws.Range("K22").Value = 750
$ws.Range("M22").Value = -455
$ws.Range("H27").Value = 1000
$ws.Range("I27").Value = 750
$ws.Range("K27").Value = 750
$ws.Range("M27").Value = -643
$ws.Range("H40").Value = 4767.8
$ws.Range("I40").Value = 4767.8
$ws.Range("K40").Value = 4767.8
$ws.Range("M40").Value = -4631.8
$ws.Range("H46").Value = 862.25
$ws.Range("I46").Value = 816.3333
$ws.Range("K46").Value = 816.3333
$ws.Range("M46").Value = -628.3333
$ws.Range("H122").Value = 6000
$ws.Range("I122").Value = 6000
$ws.Range("K122").Value = 18000
$ws.Range("M122").Value = -15550
$ws.Range("H126").Value = 4981
$ws.Range("I126").Value = 4701.1113
$ws.Range("K126").Value = 14103.3339
$ws.Range("M126").Value = -11633.3339
$ws.Range("H136").Value = 1488
$ws.Range("I136").Value = 1488
$ws.Range("K136").Value = 4464
$ws.Range("M136").Value = -1914

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 5663.3335
$ws.Range("I126").Value = 5000
$ws.Range("K126").Value = 15000
$ws.Range("M126").Value = -12530
$ws.Range("H132").Value = 2973.484
$ws.Range("I132").Value = 2049.5881
$ws.Range("J132").Value = 4095.3572
$ws.Range("K132").Value = 6148.7643
$ws.Range("L132").Value = 12286.0716
$ws.Range("M132").Value = -3618.7643
$ws.Range("N132").Value = -17346.0716
$ws.Range("H136").Value = 1635.9
$ws.Range("I136").Value = 1673.3334
$ws.Range("J136").Value = 1299
$ws.Range("K136").Value = 5020.0002
$ws.Range("L136").Value = 3897
$ws.Range("M136").Value = -2470.0002
$ws.Range("N136").Value = -8997
